$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 394: the "Numero de page" value "NA" is replaced with an empty value.
# Force text formatting first so the blank write doesn't pick up stray
# numeric/date formatting, then drop the formatting override again so the
# cell ends up plain (matching the rest of the sheet).
$ws.Range("C394").NumberFormat = "@"
$ws.Range("C394").Value = ""
$ws.Range("C394").ClearFormats()

# New rows 395-398 appended at the bottom of the table.
$ws.Range("A395:A398").NumberFormat = "@"

$ws.Range("A395").Value = "2026-01-09"
$ws.Range("B395").Value = "buse"
$ws.Range("C395").Value = 41
$ws.Range("D395").Value = 2

$ws.Range("A396").Value = "2026-01-09"
$ws.Range("B396").Value = "buse"
$ws.Range("C396").Value = 46
$ws.Range("D396").Value = 1

$ws.Range("A397").Value = "2026-01-09"
$ws.Range("B397").Value = "agriculture biologique"
$ws.Range("C397").Value = 76
$ws.Range("D397").Value = 1

$ws.Range("A398").Value = "2026-01-09"
$ws.Range("B398").Value = "agriculture biologique"
$ws.Range("C398").Value = 77
$ws.Range("D398").Value = 1

# Drop the temporary text-format override on the date column so these new
# cells end up with plain/default formatting like the rest of the sheet.
$ws.Range("A395:A398").ClearFormats()
